# Update the "pensions received" row (row 4) figures for the years
# 2015-2021 (columns E-K) with the corrected/updated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("გორი")

$ws.Range("E4").Value = 21340
$ws.Range("F4").Value = 21674
$ws.Range("G4").Value = 22026
$ws.Range("H4").Value = 22303
$ws.Range("I4").Value = 22793
$ws.Range("J4").Value = 23545
$ws.Range("K4").Value = 24002
